$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 2 }

# Row 1 holds headers; data starts on row 2 ("Förändrad" / Changed column = C).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45203
    }
}
